# Swap the deck's applied design theme from "Integral" to "Office Theme".
#
# Source diff: ppt/theme/theme1.xml and ppt/theme/theme2.xml trade places --
# theme1.xml becomes the old "Integral" theme and theme2.xml becomes the old
# "Office Theme". The slide master (used by every slide in this deck) is the
# one wired to theme2.xml, so the user-visible effect is: the deck's color
# scheme changes from the Integral palette to the stock Office palette.
#
# Reproduce that with the Design / ThemeColorScheme object model, writing the
# Office Theme's twelve standard theme colors (dk1/lt1/dk2/lt2/accent1-6/
# hlink/folHlink) over the slide master's current (Integral) theme colors, in
# COM's fixed index order.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# index -> BGR-packed RGB long (PowerPoint COM stores colors as 0x00BBGGRR),
# taken from the "Office Theme" palette: dk1, lt1, dk2, lt2, accent1..accent6,
# hlink, folHlink.
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeThemeColors[$i - 1]
}

Write-Host "Applied Office Theme color scheme to slide master."
